$d = $word.ActiveDocument

function ReplaceOnce($findText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $findText"
    }
    return $found
}

# 1. "Ator está " -> "Ator estar " (pré-condições)
ReplaceOnce "Ator está " "Ator estar "

# 2. "...para tela de semestres." -> "...para tela de Semestres." (capitaliza Semestres)
ReplaceOnce "para tela de semestres." "para tela de Semestres."

# 3. " semestre, " -> " Semestre, " (capitaliza Semestre na frase "Na tela de semestre,")
ReplaceOnce " semestre, " " Semestre, "

# 4. "semestre com status igual a “aberto”" -> "semestre aberto"
ReplaceOnce "semestre com status igual a “aberto”" "semestre aberto"

# 5. ", trabalhos-disciplina)." -> ", curso-semestre)." (primeira ocorrência, já única)
ReplaceOnce ", trabalhos-disciplina)." ", curso-semestre)."

# 6. "aluno-disciplina, professor-disciplina, trabalhos-disciplina)" -> "...curso-semestre)."
#    (segunda ocorrência; após o passo 5 esta é a única remanescente)
ReplaceOnce "aluno-disciplina, professor-disciplina, trabalhos-disciplina)" "aluno-disciplina, professor-disciplina, curso-semestre)."

# 7. Move o bookmark oculto _GoBack do parágrafo "N/A" para o final do parágrafo
#    que agora termina em "curso-semestre)."
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$endRng = $d.Content
$endRng.Find.Execute("perdendo todos os trabalhos listados e os relacionamentos (aluno-disciplina, professor-disciplina, curso-semestre).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRng)

Write-Host "done"
